# Apply "Natmi following Dr Hou advice" update to Ptn-Ptprs sheet.
# Rebuilds the full result table: the sending/target cluster set now
# includes "ECs" (previously missing), producing the complete 3x3
# (ECs/FAPs/sCs) x (ECs/FAPs/sCs) combination matrix with refreshed
# NATMI scores for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity")

for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$data = @(
        @("ECs", "Ptn", "Ptprs", "ECs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 3.090355666666666, 9.271066999999999, 0.06928583878088775, 0.06928583878088775, 0.9910080443567777, 8.919072399210998, 0.002079945264154307, 0.002079945264154307),
        @("ECs", "Ptn", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 25.17096033333333, 75.512881, 0.5643334579338453, 0.5643334579338454, 8.071764827452554, 72.64588344707299, 0.01694116321439568, 0.01694116321439569),
        @("ECs", "Ptn", "Ptprs", "sCs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 16.34167533333333, 49.025026, 0.366380703285267, 0.366380703285267, 5.240410315317555, 47.163692837858, 0.01099866613559602, 0.01099866613559602),
        @("FAPs", "Ptn", "Ptprs", "ECs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 3.090355666666666, 9.271066999999999, 0.06928583878088775, 0.06928583878088775, 17.77657358223788, 159.989162240141, 0.03730978799396602, 0.03730978799396603),
        @("FAPs", "Ptn", "Ptprs", "FAPs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 25.17096033333333, 75.512881, 0.5643334579338453, 0.5643334579338454, 144.7902690707847, 1303.112421637063, 0.3038883853307915, 0.3038883853307916),
        @("FAPs", "Ptn", "Ptprs", "sCs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 16.34167533333333, 49.025026, 0.366380703285267, 0.366380703285267, 94.00179958359976, 846.016196252398, 0.1972926445746398, 0.1972926445746399),
        @("sCs", "Ptn", "Ptprs", "ECs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 3.090355666666666, 9.271066999999999, 0.06928583878088775, 0.06928583878088775, 14.24425997096989, 128.198339738729, 0.0298961055227674, 0.0298961055227674),
        @("sCs", "Ptn", "Ptprs", "FAPs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 25.17096033333333, 75.512881, 0.5643334579338453, 0.5643334579338454, 116.0195593582608, 1044.176034224347, 0.243503909388658, 0.2435039093886581),
        @("sCs", "Ptn", "Ptprs", "sCs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 16.34167533333333, 49.025026, 0.366380703285267, 0.366380703285267, 75.32306857749577, 677.907617197462, 0.1580893925750311, 0.1580893925750311)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($col = 1; $col -le $rowVals.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowVals[$col - 1]
    }
}
